$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 1103.6666
$ws.Range("I20").Value = 1103.6666
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1103.6666
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -873.6666
$ws.Range("N20").ClearContents()
$ws.Range("H35").Value = 1103.6666
$ws.Range("I35").Value = 1103.6666
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1103.6666
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -724.6666
$ws.Range("N35").ClearContents()
$ws.Range("H76").Value = 4297
$ws.Range("I76").Value = 4256.4
$ws.Range("K76").Value = 4256.4
$ws.Range("M76").Value = -3941.4
$ws.Range("H79").Value = 4297
$ws.Range("I79").Value = 4256.4
$ws.Range("K79").Value = 4256.4
$ws.Range("M79").Value = -3164.4
$ws.Range("H125").Value = 9633
$ws.Range("I125").Value = 4999.5
$ws.Range("K125").Value = 44995.5
$ws.Range("M125").Value = -42535.5
$ws.Range("H127").Value = 3973
$ws.Range("I127").Value = 1973.5
$ws.Range("K127").Value = 5920.5
$ws.Range("M127").Value = -960.5
$ws.Range("H132").Value = 1581.5294
$ws.Range("I132").Value = 1349.0834
$ws.Range("K132").Value = 4047.2502
$ws.Range("M132").Value = -1517.2502
$ws.Range("H138").Value = 1712.7407
$ws.Range("I138").Value = 707.3333
$ws.Range("K138").Value = 2121.9999
$ws.Range("M138").Value = 3018.0001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2634.375
$ws.Range("I63").Value = 1689.6666
$ws.Range("K63").Value = 1689.6666
$ws.Range("M63").Value = -1003.6666
$ws.Range("H66").Value = 2634.375
$ws.Range("I66").Value = 1689.6666
$ws.Range("K66").Value = 8448.333000000001
$ws.Range("M66").Value = -5016.333000000001
$ws.Range("H88").Value = 2244.4443
$ws.Range("I88").Value = 751.2
$ws.Range("J88").Value = 2818.7693
$ws.Range("K88").Value = 751.2
$ws.Range("L88").Value = 2818.7693
$ws.Range("M88").Value = -345.2
$ws.Range("N88").Value = -3630.7693
$ws.Range("H91").Value = 2244.4443
$ws.Range("I91").Value = 751.2
$ws.Range("J91").Value = 2818.7693
$ws.Range("K91").Value = 751.2
$ws.Range("L91").Value = 2818.7693
$ws.Range("M91").Value = 652.8
$ws.Range("N91").Value = -5626.7693
$ws.Range("H110").Value = 3538.0667
$ws.Range("I110").Value = 614.4
$ws.Range("J110").Value = 4999.9
$ws.Range("K110").Value = 614.4
$ws.Range("L110").Value = 4999.9
$ws.Range("M110").Value = 1430.6
$ws.Range("N110").Value = -9089.9
$ws.Range("H132").Value = 2096.875
$ws.Range("I132").Value = 2096.875
$ws.Range("K132").Value = 6290.625
$ws.Range("M132").Value = -3760.625

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2187.9
$ws.Range("J86").Value = 5637.25
$ws.Range("L86").Value = 5637.25
$ws.Range("N86").Value = -7883.25
$ws.Range("H89").Value = 2187.9
$ws.Range("J89").Value = 5637.25
$ws.Range("L89").Value = 28186.25
$ws.Range("N89").Value = -39418.25
$ws.Range("H107").Value = 1330.7142
$ws.Range("I107").Value = 1330.7142
$ws.Range("K107").Value = 1330.7142
$ws.Range("M107").Value = 589.2858000000001
$ws.Range("H134").Value = 11466.111
$ws.Range("I134").Value = 13365.833
$ws.Range("K134").Value = 40097.499
$ws.Range("M134").Value = -37562.499

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 248
$ws.Range("I22").Value = 226
$ws.Range("K22").Value = 226
$ws.Range("M22").Value = 124

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 10000
$ws.Range("I104").Value = 10000
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 30000
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = -27379
$ws.Range("N104").ClearContents()
$ws.Range("H108").Value = 421.33334
$ws.Range("I108").Value = 421.33334
$ws.Range("K108").Value = 1264.00002
$ws.Range("M108").Value = 1615.99998
$ws.Range("H109").Value = 1792.6
$ws.Range("I109").Value = 1685.2
$ws.Range("K109").Value = 5055.6
$ws.Range("M109").Value = -4015.6
$ws.Range("H126").Value = 600
$ws.Range("I126").Value = 600
$ws.Range("K126").Value = 1800
$ws.Range("M126").Value = 3140
$ws.Range("H129").Value = 118.75
$ws.Range("I129").Value = 125
$ws.Range("J129").Value = 116.666664
$ws.Range("K129").Value = 375
$ws.Range("L129").Value = 349.999992
$ws.Range("M129").Value = 4625
$ws.Range("N129").Value = -10349.999992
$ws.Range("H137").Value = 3479
$ws.Range("I137").Value = 2813.6667
$ws.Range("J137").Value = 4144.3335
$ws.Range("K137").Value = 8441.000100000001
$ws.Range("L137").Value = 12433.0005
$ws.Range("M137").Value = -3341.000100000001
$ws.Range("N137").Value = -22633.0005
$ws.Range("H140").Value = 4145.8184
$ws.Range("I140").Value = 3826.375
$ws.Range("J140").Value = 4997.6665
$ws.Range("K140").Value = 11479.125
$ws.Range("L140").Value = 14992.9995
$ws.Range("M140").Value = -6299.125
$ws.Range("N140").Value = -25352.9995

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2876.3635
$ws.Range("J80").Value = 2965.1035
$ws.Range("L80").Value = 2965.1035
$ws.Range("N80").Value = -4961.1035
$ws.Range("H83").Value = 2876.3635
$ws.Range("J83").Value = 2965.1035
$ws.Range("L83").Value = 14825.5175
$ws.Range("N83").Value = -24809.5175
$ws.Range("H122").Value = 8933687
$ws.Range("I122").Value = 8933687
$ws.Range("K122").Value = 26801061
$ws.Range("M122").Value = -26798611

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 25000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 25000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 25000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -25340
$ws.Range("H55").Value = 2800
$ws.Range("I55").Value = 2700
$ws.Range("K55").Value = 2700
$ws.Range("M55").Value = -2527
$ws.Range("H82").Value = 1495.3334
$ws.Range("I82").Value = 990.8333
$ws.Range("J82").Value = 1999.8334
$ws.Range("K82").Value = 990.8333
$ws.Range("L82").Value = 1999.8334
$ws.Range("M82").Value = -629.8333
$ws.Range("N82").Value = -2721.8334
$ws.Range("H85").Value = 1495.3334
$ws.Range("I85").Value = 990.8333
$ws.Range("J85").Value = 1999.8334
$ws.Range("K85").Value = 990.8333
$ws.Range("L85").Value = 1999.8334
$ws.Range("M85").Value = 257.1667
$ws.Range("N85").Value = -4495.8334
$ws.Range("H132").Value = 12060.333
$ws.Range("J132").Value = 3086
$ws.Range("L132").Value = 9258
$ws.Range("N132").Value = -14318
$ws.Range("H136").Value = 3459.5557
$ws.Range("I136").Value = 3459.5557
$ws.Range("K136").Value = 10378.6671
$ws.Range("M136").Value = -7828.667099999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 777.4
$ws.Range("I100").Value = 679.6667
$ws.Range("K100").Value = 1359.3334
$ws.Range("M100").Value = -818.3334
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 1499.75
$ws.Range("I132").Value = 1499.6666
$ws.Range("K132").Value = 4498.9998
$ws.Range("M132").Value = -1968.9998
